$wb = $excel.ActiveWorkbook

# --- Insert new worksheet "Bank0xE" right after "Bank0xD" ----------------
# Duplicate "Bank0xD" (keeps the same column widths / header styling as the
# other Bank0xN sheets) then rename + re-populate it for the new cache data.
$template = $wb.Worksheets.Item("Bank0xD")
$template.Copy($null, $template) | Out-Null
$ws = $wb.Worksheets.Item("Bank0xD (2)")
$ws.Name = "Bank0xE"

# Header row (same as the template, set explicitly for clarity)
$ws.Range("A1").Value = "Purpose"
$ws.Range("B1").Value = "Start"
$ws.Range("C1").Value = "Size"
$ws.Range("D1").Value = "No"
$ws.Range("E1").Value = "Total Size"
$ws.Range("F1").Value = "Code Purpose"

# Row 2 - code
$ws.Range("A2").Value = "code"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 682
$ws.Range("D2").Value = 1
$ws.Range("E2").Formula = "=C2"
$ws.Range("F2").Value = "Cache Code"

# Text cells below are written in the same order the workbook's shared
# string table expects them (Logic, View, Logic Data, View Data).
$ws.Range("A3").Value = "lruCache Logic"
$ws.Range("A5").Value = "lruCache View"
$ws.Range("A4").Value = "lruCache Logic Data"
$ws.Range("A6").Value = "IruCache View Data"

# Row 3 - lruCache Logic
$ws.Range("B3").Formula = "=8191 -E3"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 1
$ws.Range("E3").Formula = "=C3 * D3"

# Row 4 - lruCache Logic Data
$ws.Range("B4").Formula = "=B3-E4"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 10
$ws.Range("E4").Formula = "=C4 * D4"

# Row 5 - lruCache View
$ws.Range("B5").Formula = "=B4-E5"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Formula = "=C5 * D5"

# Row 6 - IruCache View Data
$ws.Range("B6").Formula = "=B5-E6"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 20
$ws.Range("E6").Formula = "=C6 * D6"

# Row 12 - Total (already present via the copied template; keep as-is)
$ws.Range("A12").Value = "Total"
$ws.Range("E12").Formula = "=SUM(E2:E10)"

# View: new sheet opens with A3 selected, and becomes the active tab
$ws.Range("A3").Select() | Out-Null
$ws.Activate() | Out-Null

# --- Selection tweak on the (renumbered) BANK61 sheet ---------------------
$bank61 = $wb.Worksheets.Item("BANK61")
$bank61.Range("A6").Select() | Out-Null

# Re-activate the new sheet last so it's the one left showing/selected.
$ws.Activate() | Out-Null

Write-Output "done"
